$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM")

# Clear the L2 cell content (was "PASS")
$ws.Range("L2").ClearContents()

# Update selection / view state: top-left cell G1, active cell L2
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("L2").Select()

# Adjust workbook window width
$excel.ActiveWindow.Width = 19440
